$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21630.05437381669
$ws.Range("D2").Value = 945.00210025

$ws.Range("B3").Value = 20220.17391350002
$ws.Range("D3").Value = 865.5529498166667

$ws.Range("B4").Value = 21601.81707523336
$ws.Range("D4").Value = 952.5652421166667

$ws.Range("B5").Value = 20845.57550671669
$ws.Range("D5").Value = 895.9074810833333

$ws.Range("B6").Value = 21646.96777590003
$ws.Range("D6").Value = 958.6311651166666

$ws.Range("B7").Value = 20920.62628503336
$ws.Range("D7").Value = 902.7633239833334

$ws.Range("B8").Value = 21663.58082261669
$ws.Range("D8").Value = 954.9713586666667

$ws.Range("B9").Value = 21659.76519851669
$ws.Range("D9").Value = 919.09625145

$ws.Range("B10").Value = 20917.60760915003
$ws.Range("D10").Value = 917.0988988833333

$ws.Range("B11").Value = 21599.46352998336
$ws.Range("D11").Value = 927.3271675333333

$ws.Range("B12").Value = 20926.46264141669
$ws.Range("D12").Value = 913.3330163666667

$ws.Range("B13").Value = 20969.69194388336
$ws.Range("D13").Value = 874.8089380166667
